$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column C for rows 2 through 28
# from serial date 45472 (2024-06-29) to 45473 (2024-06-30)
for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45472) {
        $cell.Value2 = 45473
    }
}
